# Insert a new weekly price record for "Espárragos" (Vega Modelo de Temuco)
# as row 33, pushing all subsequent rows (old rows 33..125) down by one
# (new rows 34..126). This mirrors a new CSV/report row being prepended
# into the middle of the already-sorted-by-category block of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 33 downwards (old row 33 -> new row 34, ..., old row 125 -> new row 126)
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record's data
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 45274
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 300000000
$ws.Range("G33").Value = "Espárragos"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 300
$ws.Range("K33").Value = 1800
$ws.Range("L33").Value = 1800
$ws.Range("M33").Value = 1800
$ws.Range("N33").Value = "`$/kilo"
$ws.Range("O33").Value = "Región del Maule"
$ws.Range("P33").Value = 1800
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the
# rest of the "Fecha" column.
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
